$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# startDate value D2: 2024-08-01 -> 2025-04-01
$ws.Range("D2").Value = "2025-04-01"

# G2: collapse rich text run "Out-patient Specialists" into a single plain-text string
$ws.Range("G2").Value = "Scans & Diagnostic Tests/Out-patient Consultations/Out-patient Medicines/Accommodation Type/Annual Limit/Dental/Optical Benefits/Maternity (Consultations, Scans and Delivery)/Maternity Waiting Period/Complications of Pregnancy/New Born Cover"

# Column I width change (20.84 -> 44.33)
$ws.Columns.Item(9).ColumnWidth = 43.43

# Row 2 height change (62.4 -> 50.2)
$ws.Rows.Item(2).RowHeight = 50.2

# Selection change to D3
$ws.Range("D3").Select()
